# Weekly update: insert a new price-observation row for "Santina" cherries
# above the existing row 26, shifting the remaining rows (old 26-37) down
# to rows 27-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 26; everything from row 26 down (including
# the old row 37) shifts down by one, so the old row 37 becomes row 38.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new record.
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value = "Los Lagos"
$ws.Range("D26").Value = 44523
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100103
$ws.Range("H26").Value = "Frutos de hueso (carozo)"
$ws.Range("I26").Value = 100103001
$ws.Range("J26").Value = "Cereza"
$ws.Range("K26").Value = "Santina"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 400
$ws.Range("N26").Value = 27000
$ws.Range("O26").Value = 28000
$ws.Range("P26").Value = 27500
$ws.Range("Q26").Value = "$/bandeja 10 kilos"
$ws.Range("R26").Value = "Provincia de Curicó"
$ws.Range("S26").Value = 2750
$ws.Range("T26").Value = 10
